$d = $word.ActiveDocument

# --- Locate the three paragraphs around the {{generationChart}} placeholder ---
# (the paragraph with the placeholder, plus the empty paragraph right before
# it and the empty paragraph right after it) by searching for the marker
# text, rather than relying on a fixed/brittle paragraph index.
$count = $d.Paragraphs.Count
$chartIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "generationChart") {
        $chartIdx = $i
        break
    }
}

if ($chartIdx -eq -1) {
    throw "Could not locate the {{generationChart}} paragraph"
}

$pBefore = $d.Paragraphs.Item($chartIdx - 1)
$pChart  = $d.Paragraphs.Item($chartIdx)
$pAfter  = $d.Paragraphs.Item($chartIdx + 1)

# --- Remove centered alignment from all three paragraphs ---
$pBefore.Alignment = 0
$pChart.Alignment  = 0
$pAfter.Alignment  = 0

# --- Change the trailing empty paragraph's (paragraph-mark) run formatting
#     from NoProof/en-US+es-MX to Bold/en-US ---
$pAfter.Range.Select()
$sel = $word.Selection
$sel.Bold = 1
$sel.NoProofing = 0
$sel.LanguageID = "en-US"
$sel.LanguageIDFarEast = 0

Write-Output "OK"
